# Scheduled data refresh: update cached market-price / profit figures
# (currentAveragePrice*, LevePrice*, LeveProfit*) across the leve-profit
# tables on each crafting-class sheet. Values only; no formulas/formatting
# are touched.
$wb = $excel.ActiveWorkbook

# ALC row 13
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H13").Value = 22999.8
$ws.Range("J13").Value = 16249.75
$ws.Range("L13").Value = 16249.75
$ws.Range("N13").Value = -16587.75

# ALC row 15
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 3881.3062
$ws.Range("I15").Value = 3881.3062
$ws.Range("K15").Value = 11643.9186
$ws.Range("M15").Value = -11474.9186

# ALC row 70
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 11977384
$ws.Range("I70").Value = 41917492
$ws.Range("J70").Value = 1340
$ws.Range("K70").Value = 125752476
$ws.Range("L70").Value = 4020
$ws.Range("M70").Value = -125752206
$ws.Range("N70").Value = -4560

# ALC row 73
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 11977384
$ws.Range("I73").Value = 41917492
$ws.Range("J73").Value = 1340
$ws.Range("K73").Value = 125752476
$ws.Range("L73").Value = 4020
$ws.Range("M73").Value = -125751540
$ws.Range("N73").Value = -5892

# ALC row 74
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 20093740
$ws.Range("I74").Value = 28703200
$ws.Range("J74").Value = 5000
$ws.Range("K74").Value = 28703200
$ws.Range("L74").Value = 5000
$ws.Range("M74").Value = -28702264
$ws.Range("N74").Value = -6872

# ALC row 77
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H77").Value = 20093740
$ws.Range("I77").Value = 28703200
$ws.Range("J77").Value = 5000
$ws.Range("K77").Value = 143516000
$ws.Range("L77").Value = 25000
$ws.Range("M77").Value = -143511320
$ws.Range("N77").Value = -34360

# ALC row 86
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 2172.9
$ws.Range("I86").Value = 1985.5714
$ws.Range("J86").Value = 2610
$ws.Range("K86").Value = 1985.5714
$ws.Range("L86").Value = 2610
$ws.Range("M86").Value = -862.5714
$ws.Range("N86").Value = -4856

# ALC row 89
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 2172.9
$ws.Range("I89").Value = 1985.5714
$ws.Range("J89").Value = 2610
$ws.Range("K89").Value = 9927.857
$ws.Range("L89").Value = 13050
$ws.Range("M89").Value = -4311.857
$ws.Range("N89").Value = -24282

# ALC row 129
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 2461.4
$ws.Range("I129").Value = 522.5
$ws.Range("J129").Value = 10217
$ws.Range("K129").Value = 1567.5
$ws.Range("L129").Value = 30651
$ws.Range("M129").Value = 3432.5
$ws.Range("N129").Value = -40651

# ALC row 135
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 851.4545000000001
$ws.Range("I135").Value = 712.5333000000001
$ws.Range("J135").Value = 1149.1428
$ws.Range("K135").Value = 6412.7997
$ws.Range("L135").Value = 10342.2852
$ws.Range("M135").Value = -3877.7997
$ws.Range("N135").Value = -15412.2852

# ARM row 88
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 2943.3333
$ws.Range("I88").Value = 2400
$ws.Range("K88").Value = 2400
$ws.Range("M88").Value = -1994

# ARM row 91
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 2943.3333
$ws.Range("I91").Value = 2400
$ws.Range("K91").Value = 2400
$ws.Range("M91").Value = -996

# ARM row 110
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 921
$ws.Range("I110").Value = 911.125
$ws.Range("K110").Value = 911.125
$ws.Range("M110").Value = 1133.875

# ARM row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 5943.5
$ws.Range("I122").Value = 7297.3335
$ws.Range("J122").Value = 1882
$ws.Range("K122").Value = 21892.0005
$ws.Range("L122").Value = 5646
$ws.Range("M122").Value = -19442.0005
$ws.Range("N122").Value = -10546

# ARM row 123
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H123").Value = 26420
$ws.Range("J123").Value = 26420
$ws.Range("L123").Value = 26420
$ws.Range("N123").Value = -36220

# BSM row 105
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3432.3333
$ws.Range("J105").Value = 3439.6667
$ws.Range("L105").Value = 3439.6667
$ws.Range("N105").Value = -6933.6667

# CRP row 5
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 1308.75
$ws.Range("I5").Value = 450
$ws.Range("J5").Value = 2167.5
$ws.Range("K5").Value = 450
$ws.Range("L5").Value = 2167.5
$ws.Range("M5").Value = -338
$ws.Range("N5").Value = -2391.5

# CRP row 16
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1930.6
$ws.Range("I16").Value = 1720
$ws.Range("J16").Value = 2141.2
$ws.Range("K16").Value = 1720
$ws.Range("L16").Value = 2141.2
$ws.Range("M16").Value = -1433
$ws.Range("N16").Value = -2715.2

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2087.9792
$ws.Range("I31").Value = 1271.6578
$ws.Range("J31").Value = 5190
$ws.Range("K31").Value = 1271.6578
$ws.Range("L31").Value = 5190
$ws.Range("M31").Value = -976.6578
$ws.Range("N31").Value = -5780

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2087.9792
$ws.Range("I34").Value = 1271.6578
$ws.Range("J34").Value = 5190
$ws.Range("K34").Value = 1271.6578
$ws.Range("L34").Value = 5190
$ws.Range("M34").Value = -1069.6578
$ws.Range("N34").Value = -5594

# CRP row 113
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 1930.6
$ws.Range("I113").Value = 1720
$ws.Range("J113").Value = 2141.2
$ws.Range("K113").Value = 1720
$ws.Range("L113").Value = 2141.2
$ws.Range("M113").Value = 450
$ws.Range("N113").Value = -6481.2

# CUL row 5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2844.889
$ws.Range("I5").Value = 3767.3333
$ws.Range("K5").Value = 11301.9999
$ws.Range("M5").Value = -11189.9999

# CUL row 34
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 994
$ws.Range("J34").Value = 1483.3334
$ws.Range("L34").Value = 4450.0002
$ws.Range("N34").Value = -4618.0002

# CUL row 39
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 5999.75
$ws.Range("J39").Value = 5999.75
$ws.Range("L39").Value = 17999.25
$ws.Range("N39").Value = -18587.25

# CUL row 55
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 5585.615
$ws.Range("J55").Value = 7179.222
$ws.Range("L55").Value = 21537.666
$ws.Range("N55").Value = -21891.666

# CUL row 92
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 500
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 500
$ws.Range("K92").Value = 0
$ws.Range("L92").ClearContents()
$ws.Range("M92").Value = 1500
$ws.Range("N92").Value = -3996

# CUL row 113
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 688.70966
$ws.Range("I113").Value = 599.4737
$ws.Range("J113").Value = 830
$ws.Range("K113").Value = 1798.4211
$ws.Range("L113").Value = 2490
$ws.Range("M113").Value = 371.5789
$ws.Range("N113").Value = -6830

# CUL row 135
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 2844.889
$ws.Range("I135").Value = 3767.3333
$ws.Range("K135").Value = 33905.9997
$ws.Range("M135").Value = -31370.9997

# GSM row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4263.25
$ws.Range("I122").Value = 4175
$ws.Range("J122").Value = 4307.375
$ws.Range("K122").Value = 12525
$ws.Range("L122").Value = 12922.125
$ws.Range("M122").Value = -10075
$ws.Range("N122").Value = -17822.125

# GSM row 123
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 8823.789000000001
$ws.Range("J123").Value = 8823.789000000001
$ws.Range("L123").Value = 8823.789000000001
$ws.Range("N123").Value = -13723.789

# LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 112502130
$ws.Range("I122").Value = 83336170
$ws.Range("K122").Value = 250008510
$ws.Range("M122").Value = -250006060

# WVR row 81
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 139300.2
$ws.Range("I81").Value = 227800
$ws.Range("J81").Value = 50800.4
$ws.Range("K81").Value = 455600
$ws.Range("L81").Value = 101600.8
$ws.Range("M81").Value = -454539
$ws.Range("N81").Value = -103722.8

# WVR row 84
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 139300.2
$ws.Range("I84").Value = 227800
$ws.Range("J84").Value = 50800.4
$ws.Range("K84").Value = 2278000
$ws.Range("L84").Value = 508004
$ws.Range("M84").Value = -2272696
$ws.Range("N84").Value = -518612

# WVR row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 9617910
$ws.Range("I122").Value = 25001796
$ws.Range("J122").Value = 2981.25
$ws.Range("K122").Value = 75005388
$ws.Range("L122").Value = 8943.75
$ws.Range("M122").Value = -75002938
$ws.Range("N122").Value = -13843.75

# WVR row 123
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 22361.545
$ws.Range("J123").Value = 22361.545
$ws.Range("L123").Value = 22361.545
$ws.Range("N123").Value = -32161.545
